# Applies the inventory update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("فهرست محصولات")

# --- Row 2 ---
$ws.Range("F2").Value = 0
$ws.Range("I2").Value = 12763618
$ws.Range("J2").Value = 13315700
$ws.Range("K2").Value = "خیر"
$ws.Range("U2").Value = "560dffbb71d7"

# --- Row 3 ---
$ws.Range("I3").Value = 12983681
$ws.Range("J3").Value = 12983681
$ws.Range("U3").Value = "dd15d9a2e318"

# --- Row 4 ---
$ws.Range("I4").Value = 13733586
$ws.Range("J4").Value = 13733586
$ws.Range("U4").Value = "6dea89bb4663"
$ws.Range("W4").Value = 12908932

# --- Row 5 ---
$ws.Range("I5").Value = 13982425
$ws.Range("J5").Value = 13982425
$ws.Range("U5").Value = "b4d8915d08e9"
$ws.Range("W5").Value = 13142829

# --- Row 6 ---
$ws.Range("I6").Value = 12174528
$ws.Range("J6").Value = 12174528
$ws.Range("U6").Value = "b2052a9afc39"

# --- Row 7 ---
$ws.Range("I7").Value = 14881296
$ws.Range("J7").Value = 14881296
$ws.Range("U7").Value = "a578315a9058"

# --- Row 8 ---
$ws.Range("I8").Value = 8983624
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = "خیر"
$ws.Range("L8").Value = 0
$ws.Range("U8").Value = "87057c30d818"
$ws.Range("W8").Value = 8444188

# --- Row 9 ---
$ws.Range("I9").Value = 10783057
$ws.Range("J9").Value = 10783057
$ws.Range("U9").Value = "dac685ad246b"
